# "First of the trip" - add first three entries of a trip (Toronto / Niagara Falls)
# plus backfill the "Poseedor" (owner) column for a handful of earlier rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")
$tbl = $ws.ListObjects.Item("Tabla1")

# ------------------------------------------------------------------
# 1. Append the three new trip rows at the bottom of the table.
#    Doing this first (before the "Poseedor" backfill below) makes the
#    shared-strings table grow in the same order the workbook's author
#    originally typed the data.
# ------------------------------------------------------------------

$newRow = $tbl.ListRows.Add()
$r = $newRow.Range.Row
$ws.Range("A" + ($r-1) + ":H" + ($r-1)).Copy()
$ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = 45056
$ws.Range("B" + $r).Value = "Philosophers Walk"
$ws.Range("E" + $r).Value = "Philosophers Walk"
$ws.Range("F" + $r).Value = "Acuarela"
$ws.Range("G" + $r).Value = "Librea Acuarelas 2"
$ws.Range("H" + $r).Value = "Escena"

$newRow = $tbl.ListRows.Add()
$r = $newRow.Range.Row
$ws.Range("A" + ($r-1) + ":H" + ($r-1)).Copy()
$ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = 45058
$ws.Range("B" + $r).Value = "Convocation Hall"
$ws.Range("E" + $r).Value = "Convocation Hall"
$ws.Range("F" + $r).Value = "Acuarela"
$ws.Range("G" + $r).Value = "Librea Acuarelas 2"
$ws.Range("H" + $r).Value = "Arquitectura"

$newRow = $tbl.ListRows.Add()
$r = $newRow.Range.Row
$ws.Range("A" + ($r-1) + ":H" + ($r-1)).Copy()
$ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = 45059
$ws.Range("B" + $r).Value = "Niagara Falls"
$ws.Range("E" + $r).Value = "Niagara Falls"
$ws.Range("F" + $r).Value = "Acuarela"
$ws.Range("G" + $r).Value = "Librea Acuarelas 2"
$ws.Range("H" + $r).Value = "Paisaje"

# ------------------------------------------------------------------
# 2. Backfill the "Poseedor" (owner) column for a handful of older rows.
# ------------------------------------------------------------------

$ws.Range("I118").Value = "Cristina"

$ws.Range("I104").Value = "Gloria Vallejo"
$ws.Range("I105").Value = "Gloria Vallejo"
$ws.Range("I131").Value = "Gloria Vallejo"

$ws.Range("I6").Value = "?"
$ws.Range("I119").Value = "?"
$ws.Range("I126").Value = "?"
$ws.Range("I130").Value = "?"

# ------------------------------------------------------------------
# 3. Column I grew a longer value ("Gloria Vallejo"), so column 8/9
#    (shared width before) now needs to be split, matching the
#    bestFit-resize a user gets from double-clicking the column edge.
# ------------------------------------------------------------------

$ws.Columns.Item(9).ColumnWidth = 11.88671875

# ------------------------------------------------------------------
# 4. Restore the cursor/selection position as last saved by the author.
# ------------------------------------------------------------------

$ws.Activate()
$ws.Range("I7").Select()
